# Update "template_upload_karyawan.xlsx" - add 3 new "Akun System" columns
# (Username, Email, Password) with example data, matching styling of the
# existing "Mandatory" header cells (yellow fill) but with a left/right-only
# thin border, plus matching column widths and selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) -------------------------------------------------
# Order of first assignment controls the shared-string table order, so set
# them in the same order the original authors did: Username, Email, Password.
$sUsername = $ws.Range("S1")
$sUsername.Interior.Color = 65535
$sUsername.Borders.Item(7).LineStyle = 1
$sUsername.Borders.Item(10).LineStyle = 1
$sUsername.Value = "Username Akun System (Mandatory)"

$sEmail = $ws.Range("R1")
$sEmail.Interior.Color = 65535
$sEmail.Borders.Item(7).LineStyle = 1
$sEmail.Borders.Item(10).LineStyle = 1
$sEmail.Value = "Email Akun System (Mandatory)"

$sPassword = $ws.Range("T1")
$sPassword.Interior.Color = 65535
$sPassword.Borders.Item(7).LineStyle = 1
$sPassword.Borders.Item(10).LineStyle = 1
$sPassword.Value = "Password Akun System (Mandatory)"

# --- New example-data cells (row 2) -------------------------------------------
$ws.Range("R2").Value = "Ex : fathan@tcf.com"
$ws.Range("S2").Value = "Ex : FA1722"
$ws.Range("T2").Value = "Ex : passwordakun123"

# --- Column widths for the 3 new columns (R, S, T) -----------------------------
$ws.Columns.Item(18).ColumnWidth = 32.83
$ws.Columns.Item(19).ColumnWidth = 36.33
$ws.Columns.Item(20).ColumnWidth = 35.15

# --- Sheet view / selection state ----------------------------------------------
$ws.Range("S4").Select()

Write-Host "Applied Akun System upload columns"
